$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark near the approval date "6"
$d.Bookmarks("_GoBack").Delete()

# 2. Fill in the first empty row of the "Change History Log" table
#    (table 1 in the document) with the new revision entry.
$tbl = $d.Tables(1)
$row = $tbl.Rows(3)

$row.Cells(1).Range.Text = "3/21/2016"
$row.Cells(1).Range.Font.SizeBi = 12

$row.Cells(2).Range.Text = "Update after review/TFS1732/Updated text for figure 1"
$row.Cells(2).Range.Font.SizeBi = 12

$row.Cells(3).Range.Text = "Lisa Stein"
$row.Cells(3).Range.Font.SizeBi = 12

# Re-create the _GoBack bookmark around the "Change Description" cell text
$descRange = $row.Cells(2).Range
$descRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $descRange)

# 3. Update the document title text: replace the two runs
#    "Outlier Management Report " + "Data Files Requirements"
#    with "Short Duration Reporting Data File" + " Requirements",
#    keeping them as two distinct runs (and keeping the trailing
#    bookmarkEnd/line-break that immediately follow them in place).
$findRng = $d.Content.Duplicate
$found = $findRng.Find.Execute("Outlier Management Report Data Files Requirements", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Look up the numeric w:id the "_Toc321377922" bookmark (which wraps
    # this heading text) currently serializes with, so the bookmarkEnd we
    # re-emit below pairs with the correct bookmarkStart.
    $tocId = "2"
    $curXml = $d.Content.WordOpenXML
    if ($curXml -match '<w:bookmarkStart w:id="(\d+)" w:name="_Toc321377922"/>') {
        $tocId = $matches[1]
    }

    $titleStart = $findRng.Start
    # Extend one extra character past the matched text so the edit
    # range also covers the line-break run that sits right after the
    # bookmarkEnd - this keeps the bookmarkEnd anchored between the
    # new text and that break instead of being pushed out of place.
    $titleRng = $d.Range($titleStart, $findRng.End + 1)

    $titleXml = '<?xml version="1.0"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + `
        '<w:r><w:t>Short Duration Reporting Data File</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> Requirements</w:t></w:r>' + `
        "<w:bookmarkEnd w:id=`"$tocId`"/>" + `
        '<w:r w:rsidR="00CB60BC"><w:br/></w:r>' + `
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $titleRng.InsertXML($titleXml)
}
